$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-5 (columns B..M)
$values = @{
    2 = @(0.02165846824645996, 0.04358482360839844, 0.02165846824645996, 0.04358482360839844, 0.02165846824645996, 0.04358482360839844, 0.01571402549743652, 0.0327878475189209, 0.01571402549743652, 0.0327878475189209, 0.01571402549743652, 0.0327878475189209)
    3 = @(0.08838996887207032, 0.02878413200378418, 0.08838996887207032, 0.02878413200378418, 0.08838996887207032, 0.02878413200378418, 0.1059636116027832, 0.03095073699951172, 0.1059636116027832, 0.03095073699951172, 0.1059636116027832, 0.03095073699951172)
    4 = @(0.1257381439208984, 0.02428970336914062, 0.1257381439208984, 0.02428970336914062, 0.1257381439208984, 0.02428970336914062, 0.108679723739624, 0.02691059112548828, 0.108679723739624, 0.02691059112548828, 0.108679723739624, 0.02691059112548828)
    5 = @(0.03609085083007812, 0.02265701293945312, 0.03609085083007812, 0.02265701293945312, 0.03609085083007812, 0.02265701293945312, 0.02435874938964844, 0.02657623291015625, 0.02435874938964844, 0.02657623291015625, 0.02435874938964844, 0.02657623291015625)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = $i + 2  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $rowVals[$i]
    }
}

# New row 6: Ensemble
$ws.Cells.Item(6, 1).Value = "Ensemble"
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats

$row6 = @(0.2798385143280029, 0.09443140029907227, 0.2798385143280029, 0.09443140029907227, 0.2798385143280029, 0.09443140029907227, 0.3162019729614258, 0.1204162120819092, 0.3162019729614258, 0.1204162120819092, 0.3162019729614258, 0.1204162120819092)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(6, $col).Value = $row6[$i]
}
